# Update column G (K = simulated strike count) with freshly regenerated s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 2
    3 = 1
    4 = 2
    5 = 1
    6 = 1
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 2
    12 = 1
    13 = 3
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 1
    25 = 3
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 0
    33 = 1
    34 = 0
    35 = 1
    36 = 0
    37 = 2
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 0
    45 = 1
    46 = 1
    47 = 3
    48 = 3
    49 = 0
    50 = 1
    51 = 2
    53 = 1
    54 = 1
    55 = 0
    56 = 0
    57 = 1
    58 = 1
    59 = 2
    60 = 2
    61 = 0
    62 = 1
    63 = 3
    64 = 1
    66 = 2
    67 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output ("Updated " + $kValues.Count + " cells in column G (K)")
